$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "reactor_intervals": rename the "mix_culture" interval to
# "open_fermentation" and switch its reaction-model equations from equalities
# (==) to inequalities (<=), with updated coefficients - this is the core of
# "solver works with <= in reaction equations".
# ---------------------------------------------------------------------------
$wsReactor = $wb.Worksheets.Item("reactor_intervals")
$wsReactor.Range("A7").Value = "open_fermentation"
$wsReactor.Range("H7").Value = "ace <= glu*0.00213+pH*0.001+ 0.0002*casine , prop <= glu*0.0008+pH*0.00012+0.0002*casine, but <= glu*0.00343+pH*-0.007526 + 0.0002*casine"

# ---------------------------------------------------------------------------
# Sheet "connection_matrix": the same interval is referenced again in the
# header row and in column A further down - keep them in sync.
# ---------------------------------------------------------------------------
$wsConn = $wb.Worksheets.Item("connection_matrix")
$wsConn.Range("I1").Value = "open_fermentation"
$wsConn.Range("A9").Value = "open_fermentation"

# ---------------------------------------------------------------------------
# Sheet "input_output_intervals": updated input quantities/compositions.
# ---------------------------------------------------------------------------
$wsIO = $wb.Worksheets.Item("input_output_intervals")
$wsIO.Range("C2").Value = 100
$wsIO.Range("E2").Value = 0.00001
$wsIO.Range("C3").Value = 0
$wsIO.Range("E3").Value = 0.00001
$wsIO.Range("F5").Value = 1
$wsIO.Range("F6").Value = 1

# ---------------------------------------------------------------------------
# Resize/reposition the note text-box on "input_output_intervals".
# ---------------------------------------------------------------------------
$shp = $wsIO.Shapes.Item(1)
$shp.Top = 142.5
$shp.Left = 609.0
$shp.Width = 176.7
$shp.Height = 141.0

# ---------------------------------------------------------------------------
# Restore the per-sheet selections and make "input_output_intervals" the
# active tab again (it was "reactor_intervals" before).
# ---------------------------------------------------------------------------
$wsReactor.Range("F8").Select()

$wsSbml = $wb.Worksheets.Item("sbml_models")
$wsSbml.Range("A16").Select()

$wsConn.Range("C17").Select()

$wsIO.Activate()
$wsIO.Range("F15").Select()
